$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new columns before column D (shifts existing D:K data to F:M)
$ws.Range("D:E").EntireColumn.Insert()

# Step 2: Copy number formats (cell styles) from the columns that will sit beside the
# newly inserted ones so the new D/E cells pick up the same formatting as the rest of
# each row (dates in row 7/38/80, plain numbers elsewhere) without creating new styles.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G7:G102").Copy()
$ws.Range("E7:E102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Step 3: Populate the new D/E columns (newest two quarters) with the figures from the
# latest financial update, row by row. $null means "leave blank" (keeps the row's blank
# pattern), and "NA" writes the same text used elsewhere in that row.
$rowData = @(
    @(7, 43404, 43312),
    @(8, 137000, 147800),
    @(9, 112800, 123700),
    @(10, 24200, 24100),
    @(11, $null, $null),
    @(12, "NA", "NA"),
    @(13, 0, 0),
    @(14, 0, 0),
    @(15, 0, 0),
    @(16, $null, $null),
    @(17, 123100, 137500),
    @(18, 13900, 10300),
    @(19, $null, $null),
    @(20, 0, 100),
    @(21, 17300, 13700),
    @(22, 1300, 1100),
    @(23, 12500, 9200),
    @(24, 2800, 2400),
    @(25, 0, 0),
    @(26, 9800, 6800),
    @(27, 9800, 6800),
    @(28, 0, 0),
    @(29, 0, 0),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, 0, -100),
    @(33, 9800, 6800),
    @(34, 0, 0),
    @(35, 9800, 6800),
    @(38, 43404, 43312),
    @(39, $null, $null),
    @(40, $null, $null),
    @(41, 6200, 7100),
    @(42, 0, 0),
    @(43, 64700, 76500),
    @(44, 198600, 170900),
    @(45, 5300, 7200),
    @(46, 274800, 261700),
    @(47, 0, 0),
    @(48, 105500, 104200),
    @(49, 0, 0),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 3900, 3600),
    @(53, 0, 0),
    @(54, 384100, 369500),
    @(55, $null, $null),
    @(56, $null, $null),
    @(57, 0, 0),
    @(58, 111200, 110500),
    @(59, 52500, 46700),
    @(60, 163600, 157100),
    @(61, 12000, 12000),
    @(62, 8000, 8800),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 183700, 177900),
    @(67, $null, $null),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, 200500, 191500),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 200500, 191500),
    @(77, 0, 0),
    @(80, 43404, 43312),
    @(81, 9800, 6800),
    @(82, $null, $null),
    @(83, 3400, 3300),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 1400, -9500),
    @(90, $null, $null),
    @(91, 0, 0),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, 0, 0),
    @(95, $null, $null),
    @(96, 0, 0),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, 0, 0),
    @(101, 0, 0),
    @(102, 0, 0)
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $dVal = $entry[1]
    $eVal = $entry[2]
    if ($null -ne $dVal) {
        $ws.Cells.Item($r, 4).Value2 = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($r, 5).Value2 = $eVal
    }
}

# Step 4: Row 72 (Retained Earnings) - the five oldest quarters shown (F:J) become "NA"
# instead of 0, matching how this line item is reported elsewhere in the sheet.
$ws.Range("F72:J72").Value2 = "NA"

Write-Host "Edit complete"
